$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9, column A was previously stored as text "76442780"; it becomes a genuine number.
$ws.Cells.Item(9, 1).Value = 76442780

# New row 10: another payment for the same phone number 76442780.
# The phone number must stay a text value (matches the other "quote-prefixed"
# phone cell at A9 before this edit), so force text entry with a leading
# apostrophe and then drop the resulting quote-prefix style so no style
# attribute is left behind on the cell.
$ws.Cells.Item(10, 1).Value = "'76442780"
$ws.Cells.Item(10, 1).Style = "Normal"

$ws.Cells.Item(10, 2).Value = 170
$ws.Cells.Item(10, 3).Value = "Cash"
$ws.Cells.Item(10, 4).Value = "2025-08-15T09:55:46"
